$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.68539023399353
$ws.Range("B1").Value = 2.255630731582642
$ws.Range("C1").Value = 3.286489725112915
$ws.Range("D1").Value = 4.43482494354248
$ws.Range("E1").Value = 0.6409188508987427
